# Price-tracker scrape: append the latest observation (2026-02-07) as a new
# row at the bottom of the price-history table.
#   Date | Price | Discount | Incredible
#   2026-02-07 | 1977000 | 0 | 0
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$newValues = @("2026-02-07", "1977000", "0", "0")
$rng = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 4))

# Every column in this sheet (including the numeric-looking Price/Discount/
# Incredible values) is stored as plain text in the shared-string table, not
# as a real number/date. Force text ("@") number format before writing so
# Excel doesn't auto-convert "2026-02-07" into a date serial or "1977000"/"0"
# into numbers, then drop back to the default "Normal" style so no stray
# per-cell formatting is left on the new row.
$rng.NumberFormat = "@"
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($newRow, $i + 1).Value = $newValues[$i]
}
$rng.Style = "Normal"
